$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ForecastReader now maps department names to their correct monthly
# forecast figures, which changes the "Suma godzin" (B) budget values
# for a few days and ripples through the dependent "Udział w godzinach"
# (C), "Idealne godziny" (F) and "Różnica godzin" (G) columns, as well
# as the totals row (34).

    $ws.Range("B4").Value = 120000.0
    $ws.Range("C4").Value = 0.02326934264107039
    $ws.Range("F4").Value = 144.56949132887596
    $ws.Range("G4").Value = 129.53050867112407
    $ws.Range("B5").Value = 200000.0
    $ws.Range("C5").Value = 0.03878223773511732
    $ws.Range("F5").Value = 240.94915221479326
    $ws.Range("G5").Value = -16.24915221479327
    $ws.Range("B6").Value = 120000.0
    $ws.Range("C6").Value = 0.02326934264107039
    $ws.Range("F6").Value = 144.56949132887596
    $ws.Range("G6").Value = 18.930508671124045
    $ws.Range("B7").Value = 120000.0
    $ws.Range("C7").Value = 0.02326934264107039
    $ws.Range("F7").Value = 144.56949132887596
    $ws.Range("G7").Value = 73.13050867112403
    $ws.Range("C8").Value = 0.01939111886755866
    $ws.Range("F8").Value = 120.47457610739663
    $ws.Range("G8").Value = 148.22542389260337
    $ws.Range("C9").Value = 0.044599573395384916
    $ws.Range("F9").Value = 277.0915250470123
    $ws.Range("G9").Value = 2.8084749529876945
    $ws.Range("C10").Value = 0.001357378320729106
    $ws.Range("C11").Value = 0.03393445801822765
    $ws.Range("F11").Value = 210.83050818794408
    $ws.Range("G11").Value = 35.669491812055924
    $ws.Range("C12").Value = 0.03296490207484972
    $ws.Range("F12").Value = 204.8067793825743
    $ws.Range("G12").Value = 20.19322061742571
    $ws.Range("C13").Value = 0.04653868528214078
    $ws.Range("F13").Value = 289.1389826577519
    $ws.Range("G13").Value = -60.738982657751905
    $ws.Range("C14").Value = 0.001357378320729106
    $ws.Range("C15").Value = 0.05235602094240838
    $ws.Range("F15").Value = 325.2813554899709
    $ws.Range("G15").Value = -43.58135548997092
    $ws.Range("C16").Value = 0.058173356602675974
    $ws.Range("F16").Value = 361.4237283221899
    $ws.Range("G16").Value = -92.20706165552286
    $ws.Range("C17").Value = 0.002714756641458212
    $ws.Range("C18").Value = 0.03878223773511732
    $ws.Range("F18").Value = 240.94915221479326
    $ws.Range("G18").Value = 53.25084778520673
    $ws.Range("C19").Value = 0.031025790188093854
    $ws.Range("F19").Value = 192.7593217718346
    $ws.Range("G19").Value = 0.3573448948324085
    $ws.Range("C20").Value = 0.029086678301337987
    $ws.Range("F20").Value = 180.71186416109495
    $ws.Range("G20").Value = -8.711864161094951
    $ws.Range("C21").Value = 0.029086678301337987
    $ws.Range("F21").Value = 180.71186416109495
    $ws.Range("G21").Value = 18.071469172238068
    $ws.Range("C22").Value = 0.03878223773511732
    $ws.Range("F22").Value = 240.94915221479326
    $ws.Range("G22").Value = 22.45084778520672
    $ws.Range("C23").Value = 0.058173356602675974
    $ws.Range("F23").Value = 361.4237283221899
    $ws.Range("G23").Value = -121.62372832218989
    $ws.Range("C24").Value = 0.002714756641458212
    $ws.Range("C25").Value = 0.03878223773511732
    $ws.Range("F25").Value = 240.94915221479326
    $ws.Range("G25").Value = 37.30084778520674
    $ws.Range("C26").Value = 0.031025790188093854
    $ws.Range("F26").Value = 192.7593217718346
    $ws.Range("G26").Value = 28.740678228165393
    $ws.Range("C27").Value = 0.029086678301337987
    $ws.Range("F27").Value = 180.71186416109495
    $ws.Range("G27").Value = 8.288135838905049
    $ws.Range("C28").Value = 0.03296490207484972
    $ws.Range("F28").Value = 204.8067793825743
    $ws.Range("G28").Value = -3.0067793825742797
    $ws.Range("C29").Value = 0.04847779716889664
    $ws.Range("F29").Value = 301.18644026849154
    $ws.Range("G29").Value = -42.68644026849154
    $ws.Range("C30").Value = 0.058173356602675974
    $ws.Range("F30").Value = 361.4237283221899
    $ws.Range("G30").Value = -86.52372832218992
    $ws.Range("C31").Value = 0.05235602094240838
    $ws.Range("F31").Value = 325.2813554899709
    $ws.Range("G31").Value = -105.74802215663792
    $ws.Range("C32").Value = 0.04072134962187318
    $ws.Range("F32").Value = 252.99660982553291
    $ws.Range("G32").Value = -2.4132764921999126
    $ws.Range("C33").Value = 0.03878223773511732
    $ws.Range("F33").Value = 240.94915221479326
    $ws.Range("G33").Value = -13.049152214793253
    $ws.Range("B34").Value = 5157000.0
    $ws.Range("F34").Value = 6162.274567893339
